$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.04076172799836216
$ws.Range("C2").Value = 1.393451858688113
$ws.Range("D2").Value = 6.301463374471818
$ws.Range("E2").Value = 2.510271573848498
$ws.Range("F2").Value = 2.571923783630169
$ws.Range("G2").Value = 21
